# Decrease column E ("剩余") by 1 for every data row (rows 2-99),
# except row 36 which stays unchanged, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)  # Column E is the 5th column
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current - 1
    }
}
